$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ancient Nightmare 1"
$ws.Range("A3").Value = "Ancient Nightmare 2"
$ws.Range("A4").Value = "Ancient Nightmare 3"
$ws.Range("A5").Value = "Haunted Carriage 1"
$ws.Range("A6").Value = "Haunted Carriage 2"
$ws.Range("A7").Value = "Haunted Carriage 3"
$ws.Range("A8").Value = "Demon Gates 1"
$ws.Range("A9").Value = "Demon Gates 2"
$ws.Range("A10").Value = "Demon Gates 3"
$ws.Range("A11").Value = "Ancient Arena"
$ws.Range("A12").Value = "Wrathborne Invasion 1"
$ws.Range("A13").Value = "Wrathborne Invasion 2"

$ws.Range("F9").Select() | Out-Null
